$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their text representation (values like "1.002" would
# otherwise be auto-converted to numbers by Excel), then restore the default style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.107.48"
$ws.Range("E2").Value = "  -2.68%  "
$ws.Range("D3").Value = "1.716.62"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "307.32"
$ws.Range("E5").Value = "  -6.31%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "0.4700"
$ws.Range("E7").Value = "  +5.21%  "
$ws.Range("D8").Value = "0.3425"
$ws.Range("E8").Value = "  -3.85%  "
$ws.Range("D9").Value = "42.14"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("D10").Value = "0.07267"
$ws.Range("E10").Value = "  -2.38%  "
$ws.Range("E11").Value = "  -4.87%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "19.92"
$ws.Range("E13").Value = "  -4.88%  "
$ws.Range("D14").Value = "5.866"
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("D15").Value = "1.718.80"
$ws.Range("E15").Value = "  -2.99%  "
$ws.Range("D16").Value = "6.889"
$ws.Range("E16").Value = "  -4.79%  "
$ws.Range("D17").Value = "89.05"
$ws.Range("E17").Value = "  -4.63%  "
$ws.Range("D18").Value = "0.00001038"
$ws.Range("E18").Value = "  -2.10%  "
$ws.Range("D19").Value = "0.06350"
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "16.51"
$ws.Range("E21").Value = "  -3.44%  "
$ws.Range("D22").Value = "5.626"
$ws.Range("E22").Value = "  -2.65%  "
$ws.Range("D23").Value = "27.162.25"
$ws.Range("E23").Value = "  -2.67%  "
$ws.Range("D24").Value = "10.86"
$ws.Range("E24").Value = "  -3.83%  "
$ws.Range("D25").Value = "2.117"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").Value = "156.97"
$ws.Range("E26").Value = "  -3.62%  "
$ws.Range("D27").Value = "19.49"
$ws.Range("E27").Value = "  -4.26%  "
$ws.Range("D28").Value = "1.910.90"
$ws.Range("E28").Value = "  -3.27%  "
$ws.Range("D29").Value = "2.103"
$ws.Range("E29").Value = "  -2.72%  "
$ws.Range("D30").Value = "119.58"
$ws.Range("E30").Value = "  -4.42%  "
$ws.Range("D31").Value = "1.019"
$ws.Range("E31").Value = "  -7.93%  "
$ws.Range("D32").Value = "0.09161"
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("D33").Value = "3.589"
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("D34").Value = "5.321"
$ws.Range("E34").Value = "  -4.98%  "
$ws.Range("D35").Value = "0.02204"
$ws.Range("E35").Value = "  -3.79%  "
$ws.Range("D36").Value = "0.05821"
$ws.Range("E36").Value = "  -4.75%  "
$ws.Range("D37").Value = "10.98"
$ws.Range("E37").Value = "  -7.47%  "
$ws.Range("D38").Value = "0.1997"
$ws.Range("E38").Value = "  -4.94%  "
$ws.Range("D39").Value = "4.744"
$ws.Range("E39").Value = "  -4.32%  "
$ws.Range("D40").Value = "1.392"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").Value = "0.5896"
$ws.Range("E41").Value = "  -6.77%  "
$ws.Range("D42").Value = "1.121"
$ws.Range("E42").Value = "  -5.49%  "
$ws.Range("D43").Value = "7.445"
$ws.Range("E43").Value = "  -5.87%  "
$ws.Range("D44").Value = "12.60"
$ws.Range("E44").Value = "  -5.61%  "
$ws.Range("D45").Value = "0.5647"
$ws.Range("E45").Value = "  -4.24%  "
$ws.Range("D46").Value = "3.552"
$ws.Range("E46").Value = "  -5.06%  "
$ws.Range("D47").Value = "117.84"
$ws.Range("E47").Value = "  -3.70%  "
$ws.Range("D48").Value = "1.839"
$ws.Range("E48").Value = "  -5.85%  "
$ws.Range("D49").Value = "0.06649"
$ws.Range("E49").Value = "  -3.73%  "
$ws.Range("E50").Value = "  -4.51%  "
$ws.Range("D51").Value = "1.000"
$ws.Range("E51").Value = "  -0.02%  "

$ws.Range("D2:D51").Style = "Normal"
